$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C ("Förändrad") date values for rows 2-5 from 45221 to 45224
# (serial date values, equivalent to 2023-10-22 -> 2023-10-25)
$ws.Range("C2:C5").Value = 45224
